# BONDS_PROVA1.xlsx update
# - rename Sheet1 -> week1
# - rewrite week1 table to BOND/BTP/SPAIN/PORT layout (drop COUPON column + FRA/BRAZ/RUS/GER rows)
# - add a new sheet week2 with an updated snapshot of the same table, set it active

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> week1 -------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "week1"

# Clear the old table (it had columns A:E, rows 1:8) and rebuild it.
$ws1.Cells.Clear()

$ws1.Range("A1").Value = "BOND"
$ws1.Range("B1").Value = "PRICE"
$ws1.Range("C1").Value = "DURATION"
$ws1.Range("D1").Value = "CONVEXITY"

$ws1.Range("A2").Value = "BTP"
$ws1.Range("B2").Value = 129.69900000000001
$ws1.Range("C2").Value = 15.308
$ws1.Range("D2").Value = 2.782

$ws1.Range("A3").Value = "SPAIN"
$ws1.Range("B3").Value = 110.6345
$ws1.Range("C3").Value = 9.4629999999999992
$ws1.Range("D3").Value = 1.0209999999999999

$ws1.Range("A4").Value = "PORT"
$ws1.Range("B4").Value = 103.068
$ws1.Range("C4").Value = 9.7799999999999994
$ws1.Range("D4").Value = 1.0629999999999999

$ws1.Range("B3:C3").NumberFormat = "#,##0.000"

$ws1.Range("E37").Select()

# --- new sheet week2 --------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "week2"

$ws2.Range("A1").Value = "BOND"
$ws2.Range("B1").Value = "PRICE"
$ws2.Range("C1").Value = "DURATION"
$ws2.Range("D1").Value = "CONVEXITY"

$ws2.Range("A2").Value = "BTP"
$ws2.Range("B2").Value = 112.765
$ws2.Range("C2").Value = 15.29
$ws2.Range("D2").Value = 2.782

$ws2.Range("A3").Value = "SPAIN"
$ws2.Range("B3").Value = 110.697
$ws2.Range("C3").Value = 9.4459999999999997
$ws2.Range("D3").Value = 1.0169999999999999

$ws2.Range("A4").Value = "PORT"
$ws2.Range("B4").Value = 103.25700000000001
$ws2.Range("C4").Value = 9.7609999999999992
$ws2.Range("D4").Value = 1.06

$ws2.Range("C2").NumberFormat = "0.000"
$ws2.Range("C3").NumberFormat = "#,##0.000"

$ws2.Range("D11").Select()

# week2 is the sheet shown/active when the file is opened
$ws2.Activate()
